# Applies the "Updated symbol list" edit: refresh Price (col D) and Hora (col G)
# values for the crypto-price table in rows 2-51, exactly as captured by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (not auto-converted to a number), preserving
# the cells pre-existing Style so no stray formatting is introduced.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Column D ("Price") updates - only the rows whose price actually changed.
$priceUpdates = @{
    2 = "240.78"
    3 = "21.34"
    4 = "5.169"
    5 = "0.05531"
    7 = "6.330"
    8 = "0.8049"
    9 = "0.9528"
    10 = "0.1383"
    11 = "0.07308"
    12 = "0.03028"
    13 = "0.03062"
    14 = "0.09313"
    15 = "3.591"
    16 = "0.001619"
    17 = "0.04678"
    18 = "0.0005757"
    19 = "0.006433"
    20 = "0.004982"
    21 = "0.001045"
    22 = "0.0001505"
    23 = "0.0003109"
    27 = "0.1290"
    40 = "0.03827"
    41 = "0.006903"
    42 = "0.1029"
    43 = "0.003098"
    44 = "0.008283"
    45 = "0.00005961"
    46 = "0.00000000752"
    47 = "0.0005516"
    48 = "0.6844"
    49 = "0.1078"
    50 = "0.00002106"
    51 = "0.01013"
}
foreach ($row in $priceUpdates.Keys) {
    Set-TextValue $ws.Cells.Item($row, 4) $priceUpdates[$row]
}

# Column G ("Hora") - every data row moves from "22" to "23".
for ($row = 2; $row -le 51; $row++) {
    Set-TextValue $ws.Cells.Item($row, 7) "23"
}
